$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "T": add "Fixed" sub-headers above each parameter row and
# shift the "Threshold for high skill worker training" row down
# ---------------------------------------------------------------
$wsT = $wb.Worksheets.Item("T")

# Move existing row 6 (D6:E6) down to row 8, keeping both value and format
$wsT.Range("D6:E6").Copy()
$wsT.Range("D8:E8").PasteSpecial(-4122)
$wsT.Range("D6:E6").Copy()
$wsT.Range("D8:E8").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$wsT.Range("D6:E6").Clear()

# Insert "Fixed" label above each parameter row (E4 above row5, E7 above row8)
$wsT.Range("D5").Copy()
$wsT.Range("E4").PasteSpecial(-4122)
$wsT.Range("E4").Value = "Fixed"
$wsT.Range("D5").Copy()
$wsT.Range("E7").PasteSpecial(-4122)
$wsT.Range("E7").Value = "Fixed"
$excel.CutCopyMode = 0

# Widen column G (7) to fit the new "Fixed" label column
$wsT.Columns.Item(7).ColumnWidth = 35.2

$wsT.Range("E13").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "O": add a "Fixed" sub-header row above the Overhiring Limit row
# ---------------------------------------------------------------
$wsO = $wb.Worksheets.Item("O")

$wsO.Range("B6").Copy()
$wsO.Range("C5").PasteSpecial(-4122)
$wsO.Range("C5").Value = "Fixed"
$excel.CutCopyMode = 0

$wsO.Range("C5").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "S": add a "Fixed" sub-header row above the Parttime Limit row
# ---------------------------------------------------------------
$wsS = $wb.Worksheets.Item("S")

# Move existing row 5 (B5:C5) down to row 6, keeping both value and format
$wsS.Range("B5:C5").Copy()
$wsS.Range("B6:C6").PasteSpecial(-4122)
$wsS.Range("B5:C5").Copy()
$wsS.Range("B6:C6").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$wsS.Range("B5:C5").Clear()

# Insert the "Fixed" label above the relocated row
$wsS.Range("C2").Copy()
$wsS.Range("C5").PasteSpecial(-4122)
$wsS.Range("C5").Value = "Fixed"
$excel.CutCopyMode = 0

$wsS.Range("C6").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Results": just move the active selection; keep as the
# last-activated sheet so it stays the active tab in the workbook
# ---------------------------------------------------------------
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Activate()
$wsResults.Range("D2").Select() | Out-Null
